$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 134710
$ws.Range("C3").Value = 307055
$ws.Range("C4").Value = 366927
$ws.Range("C5").Value = 491801
$ws.Range("C6").Value = 584601
$ws.Range("C7").Value = 761223
$ws.Range("C8").Value = 839482
$ws.Range("C9").Value = 959225
$ws.Range("C10").Value = 1064001
$ws.Range("C11").Value = 1236345
$ws.Range("C12").Value = 1352238
$ws.Range("C13").Value = 1433065
$ws.Range("C14").Value = 1603698
$ws.Range("C15").Value = 1749100
$ws.Range("C16").Value = 2374328
$ws.Range("C17").Value = 1976184
$ws.Range("C18").Value = 2149811
$ws.Range("C19").Value = 2225078
$ws.Range("C20").Value = 2308042
$ws.Range("C21").Value = 2487657
$ws.Range("C22").Value = 2723721
$ws.Range("C23").Value = 2682239
$ws.Range("C24").Value = 3182165
$ws.Range("C25").Value = 2993998
$ws.Range("C26").Value = 3187297
$ws.Range("C27").Value = 3259998
$ws.Range("C28").Value = 3336975
$ws.Range("C29").Value = 3519584
$ws.Range("C30").Value = 3674394
$ws.Range("C31").Value = 3705612
$ws.Range("C32").Value = 4103757
$ws.Range("C33").Value = 4120863
$ws.Range("C34").Value = 4174747
$ws.Range("C35").Value = 4305182
$ws.Range("C36").Value = 4435616
$ws.Range("C37").Value = 4821359
$ws.Range("C38").Value = 4711024
$ws.Range("C39").Value = 6449433
$ws.Range("C40").Value = 5012092
$ws.Range("C41").Value = 6155636
$ws.Range("C42").Value = 5149369
$ws.Range("C43").Value = 5396552
$ws.Range("C44").Value = 5778018
$ws.Range("C45").Value = 5496623
$ws.Range("C46").Value = 5656992
$ws.Range("C47").Value = 5784005
$ws.Range("C48").Value = 6041024
$ws.Range("C49").Value = 6989131
$ws.Range("C50").Value = 6210803
$ws.Range("C51").Value = 7572022
$ws.Range("C52").Value = 6508449
$ws.Range("C53").Value = 6544373
$ws.Range("C54").Value = 6891198
$ws.Range("C55").Value = 6918141
$ws.Range("C56").Value = 7124269
$ws.Range("C57").Value = 8367457
$ws.Range("C58").Value = 8261826
$ws.Range("C59").Value = 7335958
$ws.Range("C60").Value = 7668244
$ws.Range("C61").Value = 7649427
$ws.Range("C62").Value = 8118135
$ws.Range("C63").Value = 7881643
$ws.Range("C64").Value = 9058543
$ws.Range("C65").Value = 8202810
$ws.Range("C66").Value = 8220772
$ws.Range("C67").Value = 10144354
$ws.Range("C68").Value = 10539933
$ws.Range("C69").Value = 9125257
$ws.Range("C70").Value = 10579705
$ws.Range("C71").Value = 9236447
$ws.Range("C72").Value = 9087196
$ws.Range("C73").Value = 9230459
$ws.Range("C74").Value = 9298457
$ws.Range("C75").Value = 9432740
$ws.Range("C76").Value = 11087329
$ws.Range("C77").Value = 9724399
$ws.Range("C78").Value = 9933522
$ws.Range("C79").Value = 10267946
$ws.Range("C80").Value = 10131524
$ws.Range("C81").Value = 10167447
$ws.Range("C82").Value = 10422756
$ws.Range("C83").Value = 10640859
$ws.Range("C84").Value = 10580988
$ws.Range("C85").Value = 11728808
$ws.Range("C86").Value = 11237007
$ws.Range("C87").Value = 12222320
$ws.Range("C88").Value = 11123679
$ws.Range("C89").Value = 11267799
$ws.Range("C90").Value = 12875773
$ws.Range("C91").Value = 13137925
$ws.Range("C92").Value = 13060520
$ws.Range("C93").Value = 11944773
$ws.Range("C94").Value = 12028165
$ws.Range("C95").Value = 12026027
$ws.Range("C96").Value = 12572140
$ws.Range("C97").Value = 12331799
$ws.Range("C98").Value = 14480327
$ws.Range("C99").Value = 12865510
$ws.Range("C100").Value = 12808631
$ws.Range("C101").Value = 12930086
